$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5837.5
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H137").Value = 2822.3157
$ws.Range("I137").Value = 1779
$ws.Range("J137").Value = 3581.0908
$ws.Range("K137").Value = 5337
$ws.Range("L137").Value = 10743.2724
$ws.Range("M137").Value = -2787
$ws.Range("N137").Value = -15843.2724

$ws.Range("H138").Value = 10268.23
$ws.Range("I138").Value = 1099.2
$ws.Range("J138").Value = 15998.875
$ws.Range("K138").Value = 3297.6
$ws.Range("L138").Value = 47996.625
$ws.Range("M138").Value = 1842.4
$ws.Range("N138").Value = -58276.625

$ws.Range("H141").Value = 5666.3335
$ws.Range("I141").Value = 5666.3335
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 16999.0005
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -11819.0005
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 497
$ws.Range("I5").Value = 496.66666
$ws.Range("J5").Value = 499
$ws.Range("K5").Value = 496.66666
$ws.Range("L5").Value = 499
$ws.Range("M5").Value = -384.66666
$ws.Range("N5").Value = -723

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H32").Value = 19460
$ws.Range("I32").Value = 19138.182
$ws.Range("K32").Value = 19138.182
$ws.Range("M32").Value = -18851.182

$ws.Range("H61").Value = 4810
$ws.Range("I61").Value = 4810
$ws.Range("K61").Value = 4810
$ws.Range("M61").Value = -4598

$ws.Range("H74").Value = 10251.6
$ws.Range("I74").Value = 10251.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10251.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -9377.6
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 10251.6
$ws.Range("I77").Value = 10251.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 51258
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -46890
$ws.Range("N77").ClearContents()

$ws.Range("H97").Value = 1391.9166
$ws.Range("I97").Value = 1154.909
$ws.Range("K97").Value = 1154.909
$ws.Range("M97").Value = -658.9090000000001

$ws.Range("H132").Value = 3524.6667
$ws.Range("I132").Value = 3049.6667
$ws.Range("K132").Value = 9149.000100000001
$ws.Range("M132").Value = -6619.000100000001

$ws.Range("H136").Value = 4810
$ws.Range("I136").Value = 4810
$ws.Range("K136").Value = 14430
$ws.Range("M136").Value = -11880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 497
$ws.Range("I4").Value = 496.66666
$ws.Range("J4").Value = 499
$ws.Range("K4").Value = 496.66666
$ws.Range("L4").Value = 499
$ws.Range("M4").Value = -381.66666
$ws.Range("N4").Value = -729

$ws.Range("H107").Value = 1280.25
$ws.Range("I107").Value = 1055.5
$ws.Range("J107").Value = 1505
$ws.Range("K107").Value = 1055.5
$ws.Range("L107").Value = 1505
$ws.Range("M107").Value = 864.5
$ws.Range("N107").Value = -5345

$ws.Range("H134").Value = 4000
$ws.Range("I134").Value = 4000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 12000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9465
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 394.4375
$ws.Range("I22").Value = 399.30768
$ws.Range("J22").Value = 373.33334
$ws.Range("K22").Value = 399.30768
$ws.Range("L22").Value = 373.33334
$ws.Range("M22").Value = -49.30768
$ws.Range("N22").Value = -1073.33334

$ws.Range("H107").Value = 78361.766
$ws.Range("I107").Value = 126400.5
$ws.Range("J107").Value = 1499.8
$ws.Range("K107").Value = 126400.5
$ws.Range("L107").Value = 1499.8
$ws.Range("M107").Value = -124480.5
$ws.Range("N107").Value = -5339.8

$ws.Range("H132").Value = 171583
$ws.Range("I132").Value = 502000
$ws.Range("K132").Value = 1506000
$ws.Range("M132").Value = -1503470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1390
$ws.Range("I28").Value = 1390
$ws.Range("K28").Value = 4170
$ws.Range("M28").Value = -3938

$ws.Range("H42").Value = 5980
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H60").Value = 866.6667
$ws.Range("I60").Value = 866.6667
$ws.Range("K60").Value = 2600.0001
$ws.Range("M60").Value = -2349.0001

$ws.Range("H62").Value = 5100
$ws.Range("J62").Value = 6650
$ws.Range("L62").Value = 19950
$ws.Range("N62").Value = -21322

$ws.Range("H65").Value = 5100
$ws.Range("J65").Value = 6650
$ws.Range("L65").Value = 59850
$ws.Range("N65").Value = -66714

$ws.Range("H68").Value = 1500
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4500
$ws.Range("N68").Value = -6122
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 1500
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 13500
$ws.Range("N71").Value = -21612
$ws.Range("M71").ClearContents()

$ws.Range("H74").Value = 6812.5
$ws.Range("J74").Value = 6812.5
$ws.Range("L74").Value = 20437.5
$ws.Range("N74").Value = -22559.5

$ws.Range("H77").Value = 6812.5
$ws.Range("J77").Value = 6812.5
$ws.Range("L77").Value = 61312.5
$ws.Range("N77").Value = -71920.5

$ws.Range("H100").Value = 4199.3335
$ws.Range("I100").Value = 4200
$ws.Range("J100").Value = 4199
$ws.Range("K100").Value = 12600
$ws.Range("L100").Value = 12597
$ws.Range("M100").Value = -11789
$ws.Range("N100").Value = -14219

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 349.83334
$ws.Range("I2").Value = 150
$ws.Range("K2").Value = 150
$ws.Range("M2").Value = -37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7108.7617
$ws.Range("I22").Value = 6737.3076
$ws.Range("J22").Value = 7712.375
$ws.Range("K22").Value = 6737.3076
$ws.Range("L22").Value = 7712.375
$ws.Range("M22").Value = -6442.3076
$ws.Range("N22").Value = -8302.375

$ws.Range("H27").Value = 7108.7617
$ws.Range("I27").Value = 6737.3076
$ws.Range("J27").Value = 7712.375
$ws.Range("K27").Value = 6737.3076
$ws.Range("L27").Value = 7712.375
$ws.Range("M27").Value = -6630.3076
$ws.Range("N27").Value = -7926.375

$ws.Range("H100").Value = 2978.2
$ws.Range("I100").Value = 2967
$ws.Range("K100").Value = 2967
$ws.Range("M100").Value = -2426

$ws.Range("H132").Value = 4000.875
$ws.Range("I132").Value = 2001
$ws.Range("K132").Value = 6003
$ws.Range("M132").Value = -3473

$ws.Range("H136").Value = 7456.125
$ws.Range("I136").Value = 3441.5
$ws.Range("K136").Value = 10324.5
$ws.Range("M136").Value = -7774.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1961.2
$ws.Range("I132").Value = 1201.75
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3605.25
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1075.25
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 3553.5
$ws.Range("I136").Value = 3553.5
$ws.Range("K136").Value = 10660.5
$ws.Range("M136").Value = -8110.5
